$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = "[name=`"Spokesman Czarny`"]  Forego the formality. I’m no critic of your... vernacular, 'Greatmouth Mob.' Or, perhaps, I should be calling you by your real name?`n"
$ws.Range("C18").Value = "[name=`"Spokesman Czarny`"]  For example... Mr. Mob, take your hometown, our 'City of Saplings and Artisans,' Ognisko.`n"
$ws.Range("C50").Value = "[name=`"Maria`"]  (Is that why the 'gardens' here are just 'fields' now...)`n"
$ws.Range("C57").Value = "[name=`"Maria`"]  When I think of 'Nearl,' the name my grandpa and my parents protected all their lives, I can’t imagine it was all just for the numbers and lines of the stock exchange.`n"
$ws.Range("C72").Value = "[name=`"Greatmouth Mob`"]  You wanna watch the heroes fall! You wanna watch the great stars crash to the ground! Thaaaat’s right! This season’s melee matches have a new rule in the mix, the 'Points Off' clause! You hear that, knights? Risk and reward in one!   `n"
$ws.Range("C73").Value = "[name=`"Greatmouth Mob`"]  Where’s the spice when your points never drop?! In this match, we’re following the 'Predator System!' Let’s keep it simple, sister: when someone goes up one, someone else is goin’ down one!  `n"
$ws.Range("C81").Value = "[name=`"Old Knight`"]  Nah... 'Points Off?' They ever have a 'Points Off' in the melees before? `n"
